$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing data rows (4 through 9) that no longer apply.
$ws.Range("A4:E9").EntireRow.Delete() | Out-Null

# Write the new "Glen Waverley" exposure site details into row 2.
$ws.Range("A2").Value = "Glen Waverley"
$ws.Range("B2").Value = "Village Century City  285-287 Springvale Road, Glen Waverley VIC 3150"
$ws.Range("C2").Value = "28/12/20 2:45pm-5:30pm"
$ws.Range("D2").Value = "2:45pm showing of Wonder Woman 1984 (Gold Class)"
$ws.Range("E2").Value = "old"

# Add a new row 3 repeating the location/site/period but with updated
# notes and an "Exist" flag of "new".
$ws.Range("A3").Value = "Glen Waverley"
$ws.Range("B3").Value = "Village Century City  285-287 Springvale Road, Glen Waverley VIC 3150"
$ws.Range("C3").Value = "28/12/20 2:45pm-5:30pm"
$ws.Range("D3").Value = "Case attended Gold Class screening Wonder Woman 1984"
$ws.Range("E3").Value = "new"

# Columns resize (bestFit-like) now that the text in A:D is different/shorter.
$ws.Range("A1").ColumnWidth = 11.333333333333332
$ws.Range("B1").ColumnWidth = 56.5
$ws.Range("C1").ColumnWidth = 20.833333333333336
$ws.Range("D1").ColumnWidth = 46.0

# Match the saved selection shown in the target workbook.
$ws.Range("D3").Select() | Out-Null
